# Update the Fitness values (column C) in the log sheet for run_14.
# New fitness values converge to a plateau of 7293 as generations progress.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13683
$ws.Range("C3").Value = 10909
$ws.Range("C4:C8").Value = 9821
$ws.Range("C9:C10").Value = 9368
$ws.Range("C11:C12").Value = 9019
$ws.Range("C13:C21").Value = 8087
$ws.Range("C22:C25").Value = 7775
$ws.Range("C26:C28").Value = 7750
$ws.Range("C29:C42").Value = 7295
$ws.Range("C43:C252").Value = 7293
